$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "vy" trial values (y-coordinate outcome). Mark column E
# with a 1 for every row (1-28) where C (y) is non-negative, matching the
# brute-force check for trajectories that are still valid for positive Y.
for ($r = 1; $r -le 28; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    if ($cVal -ge 0) {
        $ws.Cells.Item($r, 5).Value = 1
    }
}

# Row 29 totals up how many of the first 28 rows satisfied the check.
$ws.Range("E29").Formula = "=SUM(E1:E28)"

# Reset scroll position / selection to match the author's view state.
$ws.Range("E30").Select() | Out-Null
